$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New day's journal entry -> row 28 (right after the existing last entry,
# row 27). Copy the formatting of that last entry first so the new row
# picks up the same date-format / wrap-text styles already used in the
# sheet instead of minting brand new style records.
$ws.Range("A27:C27").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(28).RowHeight = 60

$ws.Range("A28").Value = 43159
$ws.Range("B28").Value = "J'ai fini la fonction qui permait de modifier et supprimer un article de la base de données. Je suis en train de commencer à faire celle pour ajouter un nouvel article dans la base de données. J'ai ajouté aussi des choses dans la documentation du projet dans la partie planification et celle pour le MLD"
$ws.Range("C28").Value = "3 périodes"

# Match the author's end-of-edit view state: scrolled near the bottom with
# the next (empty) row selected, ready for tomorrow's entry.
$ws.Activate()
$ws.Range("B29").Select()
